$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlinks that were previously added on A2 and A3
$ws.Range("A2:A3").Hyperlinks.Delete()

# Clear the contents of A2:C3 (the pasted / filled-in recipient data)
$ws.Range("A2:C3").ClearContents()

# Move the active selection to A3, as in the saved workbook
$ws.Range("A3").Select()
